$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = "3"
$ws.Range("F2").Value = "1"
$ws.Range("G2").Value = "1.026663333333333"
$ws.Range("H2").Value = "3.07999"
$ws.Range("I2").Value = "0.005876773047146149"
$ws.Range("J2").Value = "0.005909718423624847"
$ws.Range("K2").Value = "3"
$ws.Range("L2").Value = "1"
$ws.Range("M2").Value = "0.506715"
$ws.Range("N2").Value = "1.520145"
$ws.Range("O2").Value = "0.003122343715987576"
$ws.Range("P2").Value = "0.003132472094339857"
$ws.Range("Q2").Value = "0.52022571095"
$ws.Range("R2").Value = "4.68203139855"
$ws.Range("S2").Value = "0.00001834930539404194"
$ws.Range("T2").Value = "0.00001851202804741096"

# Row 3
$ws.Range("E3").Value = "3"
$ws.Range("F3").Value = "1"
$ws.Range("G3").Value = "1.026663333333333"
$ws.Range("H3").Value = "3.07999"
$ws.Range("I3").Value = "0.005876773047146149"
$ws.Range("J3").Value = "0.005909718423624847"
$ws.Range("M3").Value = "88.13219433333332"
$ws.Range("N3").Value = "264.396583"
$ws.Range("O3").Value = "0.5430646480820168"
$ws.Range("P3").Value = "0.5448262620252092"
$ws.Range("Q3").Value = "90.48209240824109"
$ws.Range("R3").Value = "814.3388316741699"
$ws.Range("S3").Value = "0.003191467686706305"
$ws.Range("T3").Value = "0.003219769798365038"

# Row 4
$ws.Range("E4").Value = "3"
$ws.Range("F4").Value = "1"
$ws.Range("G4").Value = "1.026663333333333"
$ws.Range("H4").Value = "3.07999"
$ws.Range("I4").Value = "0.005876773047146149"
$ws.Range("J4").Value = "0.005909718423624847"
$ws.Range("M4").Value = "1.5741895"
$ws.Range("N4").Value = "3.148379"
$ws.Range("O4").Value = "0.009700049718478087"
$ws.Range("P4").Value = "0.006487676741301404"
$ws.Range("Q4").Value = "1.616162639368333"
$ws.Range("R4").Value = "9.696975836209999"
$ws.Range("S4").Value = "0.00005700499074152961"
$ws.Range("T4").Value = "0.00003834034276459132"

# Row 5
$ws.Range("E5").Value = "3"
$ws.Range("F5").Value = "1"
$ws.Range("G5").Value = "1.026663333333333"
$ws.Range("H5").Value = "3.07999"
$ws.Range("I5").Value = "0.005876773047146149"
$ws.Range("J5").Value = "0.005909718423624847"
$ws.Range("M5").Value = "72.07364666666666"
$ws.Range("N5").Value = "216.22094"
$ws.Range("O5").Value = "0.4441129584835175"
$ws.Range("P5").Value = "0.4455535891391496"
$ws.Range("Q5").Value = "73.99537033228887"
$ws.Range("R5").Value = "665.9583329906"
$ws.Range("S5").Value = "0.002609951064304272"
$ws.Range("T5").Value = "0.002633096254447808"

# Row 6
$ws.Range("I6").Value = "0.002349989884237642"
$ws.Range("J6").Value = "0.002363164002216374"
$ws.Range("K6").Value = "3"
$ws.Range("L6").Value = "1"
$ws.Range("M6").Value = "0.506715"
$ws.Range("N6").Value = "1.520145"
$ws.Range("O6").Value = "0.003122343715987576"
$ws.Range("P6").Value = "0.003132472094339857"
$ws.Range("Q6").Value = "0.208026607195"
$ws.Range("R6").Value = "1.872239464755"
$ws.Range("S6").Value = "0.000007337476147683772"
$ws.Range("T6").Value = "0.000007402545291291284"

# Row 7
$ws.Range("I7").Value = "0.002349989884237642"
$ws.Range("J7").Value = "0.002363164002216374"
$ws.Range("M7").Value = "88.13219433333332"
$ws.Range("N7").Value = "264.396583"
$ws.Range("O7").Value = "0.5430646480820168"
$ws.Range("P7").Value = "0.5448262620252092"
$ws.Range("Q7").Value = "36.18176168420855"
$ws.Range("R7").Value = "325.635855157877"
$ws.Range("S7").Value = "0.001276196429479814"
$ws.Range("T7").Value = "0.00128751380988008"

# Row 8
$ws.Range("I8").Value = "0.002349989884237642"
$ws.Range("J8").Value = "0.002363164002216374"
$ws.Range("M8").Value = "1.5741895"
$ws.Range("N8").Value = "3.148379"
$ws.Range("O8").Value = "0.009700049718478087"
$ws.Range("P8").Value = "0.006487676741301404"
$ws.Range("Q8").Value = "0.6462672326001667"
$ws.Range("R8").Value = "3.877603395601"
$ws.Range("S8").Value = "0.00002279501871502569"
$ws.Range("T8").Value = "0.00001533144413305991"

# Row 9
$ws.Range("I9").Value = "0.002349989884237642"
$ws.Range("J9").Value = "0.002363164002216374"
$ws.Range("M9").Value = "72.07364666666666"
$ws.Range("N9").Value = "216.22094"
$ws.Range("O9").Value = "0.4441129584835175"
$ws.Range("P9").Value = "0.4455535891391496"
$ws.Range("Q9").Value = "29.58909087798444"
$ws.Range("R9").Value = "266.30181790186"
$ws.Range("S9").Value = "0.001043660959895118"
$ws.Range("T9").Value = "0.001052916202911943"

# Row 10
$ws.Range("G10").Value = "101.4397916666667"
$ws.Range("H10").Value = "304.319375"
$ws.Range("I10").Value = "0.5806563984702423"
$ws.Range("J10").Value = "0.5839115766945667"
$ws.Range("K10").Value = "3"
$ws.Range("L10").Value = "1"
$ws.Range("M10").Value = "0.506715"
$ws.Range("N10").Value = "1.520145"
$ws.Range("O10").Value = "0.003122343715987576"
$ws.Range("P10").Value = "0.003132472094339857"
$ws.Range("Q10").Value = "51.401064034375"
$ws.Range("R10").Value = "462.609576309375"
$ws.Range("S10").Value = "0.001813008856911539"
$ws.Range("T10").Value = "0.001829086719557717"

# Row 11
$ws.Range("G11").Value = "101.4397916666667"
$ws.Range("H11").Value = "304.319375"
$ws.Range("I11").Value = "0.5806563984702423"
$ws.Range("J11").Value = "0.5839115766945667"
$ws.Range("M11").Value = "88.13219433333332"
$ws.Range("N11").Value = "264.396583"
$ws.Range("O11").Value = "0.5430646480820168"
$ws.Range("P11").Value = "0.5448262620252092"
$ws.Range("Q11").Value = "8940.111432299513"
$ws.Range("R11").Value = "80461.00289069561"
$ws.Range("S11").Value = "0.3153339626918135"
$ws.Range("T11").Value = "0.3181303616837471"

# Row 12
$ws.Range("G12").Value = "101.4397916666667"
$ws.Range("H12").Value = "304.319375"
$ws.Range("I12").Value = "0.5806563984702423"
$ws.Range("J12").Value = "0.5839115766945667"
$ws.Range("M12").Value = "1.5741895"
$ws.Range("N12").Value = "3.148379"
$ws.Range("O12").Value = "0.009700049718478087"
$ws.Range("P12").Value = "0.006487676741301404"
$ws.Range("Q12").Value = "159.6854549238542"
$ws.Range("R12").Value = "958.1127295431248"
$ws.Range("S12").Value = "0.005632395934513774"
$ws.Range("T12").Value = "0.003788229555097972"

# Row 13
$ws.Range("G13").Value = "101.4397916666667"
$ws.Range("H13").Value = "304.319375"
$ws.Range("I13").Value = "0.5806563984702423"
$ws.Range("J13").Value = "0.5839115766945667"
$ws.Range("M13").Value = "72.07364666666666"
$ws.Range("N13").Value = "216.22094"
$ws.Range("O13").Value = "0.4441129584835175"
$ws.Range("P13").Value = "0.4455535891391496"
$ws.Range("Q13").Value = "7311.13570252361"
$ws.Range("R13").Value = "65800.22132271249"
$ws.Range("S13").Value = "0.2578770309870035"
$ws.Range("T13").Value = "0.260163898736164"

# Row 14
$ws.Range("G14").Value = "2.9217165"
$ws.Range("H14").Value = "5.843433"
$ws.Range("I14").Value = "0.01672433817506114"
$ws.Range("J14").Value = "0.01121206356427047"
$ws.Range("K14").Value = "3"
$ws.Range("L14").Value = "1"
$ws.Range("M14").Value = "0.506715"
$ws.Range("N14").Value = "1.520145"
$ws.Range("O14").Value = "0.003122343715987576"
$ws.Range("P14").Value = "0.003132472094339857"
$ws.Range("Q14").Value = "1.4804775762975"
$ws.Range("R14").Value = "8.882865457785"
$ws.Range("S14").Value = "0.00005221913220495327"
$ws.Range("T14").Value = "0.00003512147623504193"

# Row 15
$ws.Range("G15").Value = "2.9217165"
$ws.Range("H15").Value = "5.843433"
$ws.Range("I15").Value = "0.01672433817506114"
$ws.Range("J15").Value = "0.01121206356427047"
$ws.Range("M15").Value = "88.13219433333332"
$ws.Range("N15").Value = "264.396583"
$ws.Range("O15").Value = "0.5430646480820168"
$ws.Range("P15").Value = "0.5448262620252092"
$ws.Range("Q15").Value = "257.4972863649065"
$ws.Range("R15").Value = "1544.983718189439"
$ws.Range("S15").Value = "0.009082396825444216"
$ws.Range("T15").Value = "0.006108626681310525"

# Row 16
$ws.Range("G16").Value = "2.9217165"
$ws.Range("H16").Value = "5.843433"
$ws.Range("I16").Value = "0.01672433817506114"
$ws.Range("J16").Value = "0.01121206356427047"
$ws.Range("M16").Value = "1.5741895"
$ws.Range("N16").Value = "3.148379"
$ws.Range("O16").Value = "0.009700049718478087"
$ws.Range("P16").Value = "0.006487676741301404"
$ws.Range("Q16").Value = "4.59933543627675"
$ws.Range("R16").Value = "18.397341745107"
$ws.Range("S16").Value = "0.0001622269118067341"
$ws.Range("T16").Value = "0.00007274024400791046"

# Row 17
$ws.Range("G17").Value = "2.9217165"
$ws.Range("H17").Value = "5.843433"
$ws.Range("I17").Value = "0.01672433817506114"
$ws.Range("J17").Value = "0.01121206356427047"
$ws.Range("M17").Value = "72.07364666666666"
$ws.Range("N17").Value = "216.22094"
$ws.Range("O17").Value = "0.4441129584835175"
$ws.Range("P17").Value = "0.4455535891391496"
$ws.Range("Q17").Value = "210.57876268117"
$ws.Range("R17").Value = "1263.47257608702"
$ws.Range("S17").Value = "0.007427495305605233"
$ws.Range("T17").Value = "0.004995575162716995"

# Row 18
$ws.Range("G18").Value = "68.89977133333333"
$ws.Range("H18").Value = "206.699314"
$ws.Range("I18").Value = "0.3943925004233126"
$ws.Range("J18").Value = "0.3966034773153216"
$ws.Range("K18").Value = "3"
$ws.Range("L18").Value = "1"
$ws.Range("M18").Value = "0.506715"
$ws.Range("N18").Value = "1.520145"
$ws.Range("O18").Value = "0.003122343715987576"
$ws.Range("P18").Value = "0.003132472094339857"
$ws.Range("Q18").Value = "34.91254763117"
$ws.Range("R18").Value = "314.21292868053"
$ws.Range("S18").Value = "0.001231428945329358"
$ws.Range("T18").Value = "0.001242349325208395"

# Row 19
$ws.Range("G19").Value = "68.89977133333333"
$ws.Range("H19").Value = "206.699314"
$ws.Range("I19").Value = "0.3943925004233126"
$ws.Range("J19").Value = "0.3966034773153216"
$ws.Range("M19").Value = "88.13219433333332"
$ws.Range("N19").Value = "264.396583"
$ws.Range("O19").Value = "0.5430646480820168"
$ws.Range("P19").Value = "0.5448262620252092"
$ws.Range("Q19").Value = "6072.288036671562"
$ws.Range("R19").Value = "54650.59233004406"
$ws.Range("S19").Value = "0.2141806244485729"
$ws.Range("T19").Value = "0.2160799900519065"

# Row 20
$ws.Range("G20").Value = "68.89977133333333"
$ws.Range("H20").Value = "206.699314"
$ws.Range("I20").Value = "0.3943925004233126"
$ws.Range("J20").Value = "0.3966034773153216"
$ws.Range("M20").Value = "1.5741895"
$ws.Range("N20").Value = "3.148379"
$ws.Range("O20").Value = "0.009700049718478087"
$ws.Range("P20").Value = "0.006487676741301404"
$ws.Range("Q20").Value = "108.4612965853343"
$ws.Range("R20").Value = "650.767779512006"
$ws.Range("S20").Value = "0.003825626862701023"
$ws.Range("T20").Value = "0.002573035155297871"

# Row 21
$ws.Range("G21").Value = "68.89977133333333"
$ws.Range("H21").Value = "206.699314"
$ws.Range("I21").Value = "0.3943925004233126"
$ws.Range("J21").Value = "0.3966034773153216"
$ws.Range("M21").Value = "72.07364666666666"
$ws.Range("N21").Value = "216.22094"
$ws.Range("O21").Value = "0.4441129584835175"
$ws.Range("P21").Value = "0.4455535891391496"
$ws.Range("Q21").Value = "4965.857774492795"
$ws.Range("R21").Value = "44692.71997043516"
$ws.Range("S21").Value = "0.1751548201667093"
$ws.Range("T21").Value = "0.1767081027829088"

